# Added BR - Count Vectorization results for Lucene, Thunderbird, Ubuntu
#
# Fill in the F1 / Precision / Recall results for the "CV + tfidf"
# (Count Vectorization) configuration rows on the Lucene sheet.

$wb = $excel.ActiveWorkbook

$lucene = $wb.Worksheets.Item("Lucene")

$lucene.Range("C9").Value  = "0 0.54 0.24 0.01 0.54 "
$lucene.Range("D9").Value  = "0 0.78 0.84 0.5 1"
$lucene.Range("E9").Value  = "0 0.41 0.14 0.01 0.37"

$lucene.Range("C10").Value = "0.02 0.33 0.07 0 0.36"
$lucene.Range("D10").Value = "0.6 0.8 0.78 0 0.9"
$lucene.Range("E10").Value = "0.01 0.21 0.04 0 0.22"

$lucene.Range("C11").Value = "0.06 0.61 0.49 0.13 0.71"
$lucene.Range("D11").Value = "0.25 0.7 0.64 0.35 0.92"
$lucene.Range("E11").Value = "0.03 0.53 0.39 0.08 0.57"

$lucene.Range("C12").Value = "0.19 0.51 0.45 0.29 0.65"
$lucene.Range("D12").Value = "0.24 0.55 0.43 0.34 0.65 "
$lucene.Range("E12").Value = "0.16 0.48 0.48 0.26 0.65"

$lucene.Range("C13").Value = "0.02 0.44 0.34 0.04 0.67"
$lucene.Range("D13").Value = "1 0.72 0.58 0.75 0.94"
$lucene.Range("E13").Value = "0.01 0.31 0.24 0.02 0.53"

# Restore the workbook's navigation/selection state: Ubuntu was the
# active tab before this edit; afterwards Lucene (where the new data was
# entered) is the active tab, with the cursor left at C15. Ubuntu's
# cursor ends up at C17, Thunderbird is untouched.

$ubuntu = $wb.Worksheets.Item("Ubuntu")
$ubuntu.Activate()
$ubuntu.Range("C17").Select()

$lucene.Activate()
$lucene.Range("C15").Select()
